$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timelog entries get their date-column formatting (numeric date format +
# wrap text) copied from the existing A2 entry so they share the same style
# instead of minting a duplicate one.

# Row 4: 11/10/2025, 1 hour - finalizing powerpoint / articles / planning
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Value = 45971
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "Finalizing the powerpoint, finding other articles, and other planing steps"
$ws.Rows.Item(4).RowHeight = 43.2

# Row 5: 11/17/2025, 2 hours - cleaning portuguese data 90%
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Value = 45978
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "beguin cleaning portuguese data 90%"
$ws.Rows.Item(5).RowHeight = 28.8

# Row 6: 11/18/2025, 1.75 hours - finished portuguese data, began papiamento
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("A6").Value = 45979
$ws.Range("B6").Value = 1.75
$ws.Range("C6").Value = "finished cleaning the portguese data completely and beguin cleaning papiamento"
$ws.Rows.Item(6).RowHeight = 57.6

# Column A widened slightly to fit the new dates; selection moved below the
# last entry, matching where the author left off.
$ws.Columns.Item(1).ColumnWidth = 9.71
$ws.Range("B7").Select() | Out-Null

Write-Output "done"
